# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
#
# The underlying match-result dataset had several pairs of rows whose
# entire record (match id, teams, odds, etc.) were swapped - only the
# leading sequential "id" column (column A) stays put for each physical
# row. This mirrors the shared-strings reorder + row-content swap seen
# in the authoritative OOXML diff for this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($row1, $row2, $firstCol, $lastCol) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Columns B (id/matchid) through AC (29) hold the match data; column A
# (the running row index 0,1,2,...) must stay exactly where it is.
$firstCol = 2   # B
$lastCol  = 29  # AC

Swap-RowRange 2  3  $firstCol $lastCol
Swap-RowRange 11 13 $firstCol $lastCol
Swap-RowRange 14 15 $firstCol $lastCol
Swap-RowRange 40 41 $firstCol $lastCol
Swap-RowRange 46 47 $firstCol $lastCol
